$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B10").Value = 1000
$ws.Range("E10").Value = "j1.jos"
$ws.Range("F10").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("B11").Value = 724
$ws.Range("E11").Value = "j1.lat"
$ws.Range("F11").Value = 1
$ws.Range("B12").Value = 724
$ws.Range("E12").Value = "j1.lat"
$ws.Range("F12").Value = 1
$ws.Range("B13").Value = 962
$ws.Range("C13").Value = 100
$ws.Range("E13").Value = "j1.leg1"
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1
$ws.Range("B14").Value = 962
$ws.Range("C14").Value = 100
$ws.Range("E14").Value = "j1.leg2"
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("B15").Value = 964
$ws.Range("C15").Value = 100
$ws.Range("E15").Value = "j1.leg"
$ws.Range("F15").Value = 1
$ws.Range("B16").Value = 600
$ws.Range("C16").Value = 500
$ws.Range("E16").Value = "t1.jos"
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = ""
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 1
$ws.Range("B17").Value = 1952
$ws.Range("C17").Value = 550
$ws.Range("E17").Value = "t1.lat"
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = ""
$ws.Range("H17").Value = 1
$ws.Range("B18").Value = 1952
$ws.Range("C18").Value = 550
$ws.Range("E18").Value = "t1.lat"
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = 1
$ws.Range("I18").Value = ""
$ws.Range("B19").Value = 564
$ws.Range("C19").Value = 499
$ws.Range("E19").Value = "t1.sus"
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = ""
$ws.Range("I19").Value = ""
$ws.Range("A20").Value = "1"
$ws.Range("B20").Value = 564
$ws.Range("C20").Value = 500
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = "t1.sep.h"
$ws.Range("F20").Value = 1
$ws.Range("A21").Value = "1"
$ws.Range("B21").Value = 564
$ws.Range("C21").Value = 500
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = "t1.sep.h"
$ws.Range("F21").Value = 1
$ws.Range("A22").Value = "1"
$ws.Range("B22").Value = 564
$ws.Range("C22").Value = 500
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = "t1.sep.h"
$ws.Range("F22").Value = 1